$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "sports_club_coed"
$ws.Range("C2").Value = "Archery"

# Row 3
$ws.Range("B3").Value = "sports_club_coed"
$ws.Range("C3").Value = "Archery"

# Row 4
$ws.Range("B4").Value = "sports_club_boys"
$ws.Range("C4").Value = "Archery"

# Row 5
$ws.Range("B5").Value = "sports_club_girls"
$ws.Range("C5").Value = "Archery"

# Row 6
$ws.Range("B6").Value = "sports_club_coed"
$ws.Range("C6").Value = "Archery"

# Row 7
$ws.Range("B7").Value = "sports_club_coed"
$ws.Range("C7").Value = "Archery"

# Row 8
$ws.Range("B8").Value = "sports_club_coed"
$ws.Range("C8").Value = "Archery"

# Row 9
$ws.Range("B9").Value = "sports_club_coed"
$ws.Range("C9").Value = "Archery"

# Row 10
$ws.Range("B10").Value = "sports_club_coed"
$ws.Range("C10").Value = "Archery"

$wb.Save()
